# Update cryptocurrency price/volume snapshot to the latest scraped values.
# (GitHub Actions scheduled refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '22.386.22'
$ws.Range("E2").Value = '  -0.02%  '

# Row 3
$ws.Range("D3").Value = '1.572.13'
$ws.Range("E3").Value = '  +0.18%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '''1.002'
$ws.Range("E5").Value = '  +0.05%  '

# Row 6
$ws.Range("D6").Value = '''290.82'
$ws.Range("E6").Value = '  -0.03%  '

# Row 7
$ws.Range("D7").Value = '''0.3759'
$ws.Range("E7").Value = '  +2.07%  '

# Row 8
$ws.Range("D8").Value = '''49.96'
$ws.Range("E8").Value = '  +0.83%  '

# Row 9
$ws.Range("D9").Value = '''0.3421'
$ws.Range("E9").Value = '  +0.67%  '

# Row 10
$ws.Range("D10").Value = '''0.07636'

# Row 11
$ws.Range("E11").Value = '  -1.99%  '

# Row 12
$ws.Range("E12").Value = '  +0.07%  '

# Row 13
$ws.Range("E13").Value = '  -0.07%  '

# Row 14
$ws.Range("D14").Value = '''6.022'

# Row 15
$ws.Range("D15").Value = '''6.942'
$ws.Range("E15").Value = '  +0.55%  '

# Row 16
$ws.Range("D16").Value = '1.573.94'
$ws.Range("E16").Value = '  -0.67%  '

# Row 17
$ws.Range("D17").Value = '''0.00001131'
$ws.Range("E17").Value = '  -0.64%  '

# Row 18
$ws.Range("D18").Value = '''90.10'
$ws.Range("E18").Value = '  +1.11%  '

# Row 19
$ws.Range("D19").Value = '''0.06739'
$ws.Range("E19").Value = '  -0.61%  '

# Row 20
$ws.Range("E20").Value = '  -0.01%  '

# Row 21
$ws.Range("D21").Value = '''16.78'
$ws.Range("E21").Value = '  +1.66%  '

# Row 22
$ws.Range("D22").Value = '''6.201'
$ws.Range("E22").Value = '  -0.54%  '

# Row 23
$ws.Range("D23").Value = '''11.99'
$ws.Range("E23").Value = '  -0.38%  '

# Row 24
$ws.Range("D24").Value = '22.378.64'
$ws.Range("E24").Value = '  -0.14%  '

# Row 25
$ws.Range("D25").Value = '''2.399'
$ws.Range("E25").Value = '  +0.51%  '

# Row 26
$ws.Range("D26").Value = '''2.684'
$ws.Range("E26").Value = '  -10.04%  '

# Row 27
$ws.Range("D27").Value = '''20.21'
$ws.Range("E27").Value = '  +1.68%  '

# Row 28
$ws.Range("D28").Value = '''147.07'
$ws.Range("E28").Value = '  +1.01%  '

# Row 29
$ws.Range("D29").Value = '''5.023'

# Row 30
$ws.Range("D30").Value = '''126.12'
$ws.Range("E30").Value = '  +0.51%  '

# Row 31
$ws.Range("D31").Value = '1.746.17'
$ws.Range("E31").Value = '  -0.54%  '

# Row 32
$ws.Range("E32").Value = '  -1.70%  '

# Row 33
$ws.Range("D33").Value = '''2.006'
$ws.Range("E33").Value = '  +0.56%  '

# Row 34
$ws.Range("D34").Value = '''0.9822'
$ws.Range("E34").Value = '  -6.04%  '

# Row 35
$ws.Range("D35").Value = '''9.881'
$ws.Range("E35").Value = '  -4.17%  '

# Row 36
$ws.Range("D36").Value = '''0.08533'
$ws.Range("E36").Value = '  +0.89%  '

# Row 37
$ws.Range("D37").Value = '''0.02541'
$ws.Range("E37").Value = '  -0.20%  '

# Row 38
$ws.Range("D38").Value = '''0.2317'
$ws.Range("E38").Value = '  -0.43%  '

# Row 39
$ws.Range("D39").Value = '''1.347'
$ws.Range("E39").Value = '  +7.97%  '

# Row 40
$ws.Range("D40").Value = '''0.06552'
$ws.Range("E40").Value = '  -0.08%  '

# Row 41
$ws.Range("D41").Value = '''5.412'
$ws.Range("E41").Value = '  -2.17%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.6386'
$ws.Range("E42").Value = '  +0.27%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '''11.43'
$ws.Range("E43").Value = '  -3.27%  '

# Row 44
$ws.Range("D44").Value = '''1.002'

# Row 45
$ws.Range("D45").Value = '''14.07'
$ws.Range("E45").Value = '  -1.84%  '

# Row 46
$ws.Range("D46").Value = '''3.778'
$ws.Range("E46").Value = '  +0.13%  '

# Row 47
$ws.Range("D47").Value = '''0.5980'
$ws.Range("E47").Value = '  -0.20%  '

# Row 48
$ws.Range("D48").Value = '''1.299'
$ws.Range("E48").Value = '  +2.97%  '

# Row 49
$ws.Range("D49").Value = '''2.083'
$ws.Range("E49").Value = '  -2.50%  '

# Row 50
$ws.Range("D50").Value = '''125.20'
$ws.Range("E50").Value = '  +1.36%  '

# Row 51
$ws.Range("E51").Value = '  +0.44%  '
